$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (was the old row 16 / A=112178654 data -> now becomes old row18 data)
$ws.Range("A16").Value = 112178652
$ws.Range("B16").Value = 90678
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4366
$ws.Range("F16").Value = "Skarp dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum peckii"
$ws.Range("H16").Value = "Banker"
$ws.Range("Q16").Value = 618476.2382824289
$ws.Range("R16").Value = 6905001.69355389

# Row 17 (becomes old row16 data)
$ws.Range("A17").Value = 112178654
$ws.Range("B17").Value = 89686
$ws.Range("E17").Value = 658
$ws.Range("F17").Value = "Rosenticka"
$ws.Range("G17").Value = "Rhodofomes roseus"
$ws.Range("H17").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 618387.2260358589
$ws.Range("R17").Value = 6904851.227267566

# Row 18 (becomes old row17 data)
$ws.Range("A18").Value = 112178651
$ws.Range("B18").Value = 86223
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 4412
$ws.Range("F18").Value = "Äggvaxskivling"
$ws.Range("G18").Value = "Hygrophorus karstenii"
$ws.Range("H18").Value = "Sacc. & Cub."
$ws.Range("Q18").Value = 618387.9774688096
$ws.Range("R18").Value = 6904949.162718941
